$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attempt 2")

# Row 5: change label from Pattern to Percentage, and change B5:E5 from text to numbers
$ws.Range("A5").Value = "Percentage"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 0

# Rows 26-29: clear the little lookup table, leaving only empty styled A26/B26
$ws.Range("A26").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("C26").Value = ""
$ws.Rows("27:29").Delete()

